# converter jin to kg unit
# Convert the ingredient weight values in column N ("重量(公斤)" / Weight in kg)
# from "jin" (斤) denominated text to kg-denominated decimal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$weights = @{
    "N2"  = "15.0"
    "N3"  = "5.3999999999999995"
    "N4"  = "1.2"
    "N5"  = "6.0"
    "N6"  = "1.7999999999999998"
    "N7"  = "0.0"
    "N8"  = "0.0"
    "N9"  = "0.0"
    "N10" = "0.0"
    "N11" = "0.6"
    "N12" = "0.0"
}

foreach ($addr in $weights.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $weights[$addr]
}
